$wb = $excel.ActiveWorkbook

# Replace "Ready for handoff" with "In Translation" wherever it occurs
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# Narrow the "zh-cn" / "de-de" status columns now that the status text is shorter
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
